{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts,\n// large numbers) in bold + dark slate color (#2C3E50) across the resume.\n//\n// For each target paragraph we search for the specific metric substrings\n// and apply bold + font color to just that sub-range, which causes Word\n// to split the paragraph's run at the match boundaries (mirroring the\n// exact <w:r> splits seen in the target OOXML diff).\n\nconst COLOR = \"#2C3E50\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Highlights the FIRST occurrence of `text` inside `paragraph`.\nasync function highlightFirst(paragraph, text) {\n  const results = paragraph.search(text, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Highlight target not found: \"${text}\"`);\n  }\n  const target = results.items[0];\n  target.font.bold = true;\n  target.font.color = COLOR;\n}\n\n// [paragraph index, [metric substrings to bold+color, in left-to-right order]]\nconst plan = [\n  // \"Discovered systematic race coding errors ... from 23% to 64%\"\n  [9, [\"23%\", \"64%\"]],\n  // \"Utilized advanced sampling methods ... \u00b14.2% to \u00b12.1% ... 71% to 87% ...\"\n  [11, [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"]],\n  // \"Trigonometric algorithm ... reduced mapping costs by 73.5% ... $4.7M ...\"\n  [12, [\"73.5%\", \"$4.7M\"]],\n  // \"Built real-time FEC analysis systems ... valued over $2 trillion\"\n  [13, [\"$2\"]],\n  // \"Modernized legacy ETL processes ... reducing processing time by 57%\"\n  [18, [\"57%\"]],\n  // \"Platform impact: Built redistricting system serving 12,847 analysts ...\"\n  [49, [\"12,847\"]],\n  // \"Revenue generation: Delivered $4.9M additional revenue ...\"\n  [51, [\"$4.9M\"]],\n  // \"23% conversion rate improvement\"\n  [52, [\"23%\"]],\n];\n\nfor (const [idx, terms] of plan) {\n  const paragraph = paragraphs.items[idx];\n  for (const term of terms) {\n    await highlightFirst(paragraph, term);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts,\n# large numbers) in bold + dark slate color (#2C3E50) across the resume.\n#\n# For each target paragraph we use Find.Execute to locate each metric\n# substring (in left-to-right order) and bold + color just that matched\n# range. Word stores this as a run split at the match boundaries,\n# mirroring the exact <w:r> splits in the target OOXML.\n\n$d = $word.ActiveDocument\n\nfunction Highlight-Metrics($paragraphIndex, $terms) {\n    $p = $d.Paragraphs.Item($paragraphIndex)\n    $paraEnd = $p.Range.End\n    $rng = $p.Range\n\n    foreach ($text in $terms) {\n        $f = $rng.Find\n        $f.Text = $text\n        $f.MatchCase = $true\n        $f.Forward = $true\n        $f.Wrap = 0\n        $found = $f.Execute()\n        if (-not $found) {\n            throw \"Highlight target not found: '$text' (paragraph $paragraphIndex)\"\n        }\n\n        # Find.Execute collapses $rng to the matched text - format it.\n        $rng.Font.Bold = $true\n        $rng.Font.Color = '#2C3E50'\n\n        # Re-expand the range past this match through the paragraph end so\n        # the next Find (if any) can locate the following occurrence.\n        $matchEnd = $rng.End\n        $rng.Start = $matchEnd\n        $rng.End = $paraEnd\n    }\n}\n\n# Paragraphs.Item is 1-based.\n# \"Discovered systematic race coding errors ... from 23% to 64%\"\nHighlight-Metrics 10 @('23%', '64%')\n# \"Utilized advanced sampling methods ... \u00b14.2% to \u00b12.1% ... 71% to 87% ...\"\nHighlight-Metrics 12 @('\u00b14.2%', '\u00b12.1%', '71%', '87%')\n# \"Trigonometric algorithm ... reduced mapping costs by 73.5% ... $4.7M ...\"\nHighlight-Metrics 13 @('73.5%', '$4.7M')\n# \"Built real-time FEC analysis systems ... valued over $2 trillion\"\nHighlight-Metrics 14 @('$2')\n# \"Modernized legacy ETL processes ... reducing processing time by 57%\"\nHighlight-Metrics 19 @('57%')\n# \"Platform impact: Built redistricting system serving 12,847 analysts ...\"\nHighlight-Metrics 50 @('12,847')\n# \"Revenue generation: Delivered $4.9M additional revenue ...\"\nHighlight-Metrics 52 @('$4.9M')\n# \"23% conversion rate improvement\"\nHighlight-Metrics 53 @('23%')\n"}
